$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 139 (shifts existing rows 139.. down by one)
$ws.Rows(139).Insert()

# Populate the new row with the New Zealand Election Study entry
$ws.Range("A139").Value = "New Zealand Election Study"
$ws.Range("B139").Value = "citizens"
$ws.Range("C139").Value = "http://www.nzes.org/exec/show/data"
$ws.Range("D139").Value = "Political attitudes, political participation, socio-demographic characteristics"
$ws.Range("E139").Value = "New Zealand"
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 1
$ws.Range("K139").Value = 1990
$ws.Range("L139").Value = 2014
$ws.Range("M139").Value = "online"
$ws.Range("N139").Value = "no"
$ws.Range("O139").Value = 1
$ws.Range("V139").Value = "http://www.nzes.org/data/NZES2014GeneralReleaseApril16.sav.zip"
$ws.Range("AB139").Value = 20180307

# Remove stray formatted-but-empty cells that the row insert copied down
# from row 138 (P139, T139 had the Hyperlink style carried over)
$ws.Range("P139").Clear()
$ws.Range("T139").Clear()

# Add hyperlinks for the link cells, then restore their Hyperlink style
# (Hyperlinks.Add applies its own default style otherwise)
$ws.Hyperlinks.Add($ws.Range("C139"), "http://www.nzes.org/exec/show/data")
$ws.Range("C139").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("V139"), "http://www.nzes.org/data/NZES2014GeneralReleaseApril16.sav.zip")
$ws.Range("V139").Style = "Hyperlink"

# Update sheet view: drop the stale topLeftCell/selection from before the
# insert, matching the post-edit scroll position
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A221").Select()
